$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-10 07:02:50"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
